$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Determine the used range to find the last row with data.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 443 }

# Column C holds the "Förändrad" (Changed) date. Every data row (2..lastRow)
# currently stores the serial date 45188 (2023-09-19) and must be bumped to
# 45189 (2023-09-20), keeping the existing date style/format untouched.
$range = $ws.Range("C2:C$lastRow")
$range.Value = 45189
